$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 426.87097
$ws.Range("I28").Value = 421.33334
$ws.Range("J28").Value = 445.85715
$ws.Range("K28").Value = 421.33334
$ws.Range("L28").Value = 445.85715
$ws.Range("M28").Value = 63.66665999999998
$ws.Range("N28").Value = -1415.85715
$ws.Range("H100").Value = 2025.0714
$ws.Range("I100").Value = 1852.4445
$ws.Range("K100").Value = 1852.4445
$ws.Range("M100").Value = -1311.4445
$ws.Range("H107").Value = 37037530
$ws.Range("I107").Value = 66666810
$ws.Range("J107").Value = 924.75
$ws.Range("K107").Value = 66666810
$ws.Range("L107").Value = 924.75
$ws.Range("M107").Value = -66664890
$ws.Range("N107").Value = -4764.75
$ws.Range("H111").Value = 3971023.8
$ws.Range("I111").Value = 5557472.5
$ws.Range("K111").Value = 16672417.5
$ws.Range("M111").Value = -16669350.5
$ws.Range("H137").Value = 102496.555
$ws.Range("I137").Value = 163470.1
$ws.Range("K137").Value = 490410.3
$ws.Range("M137").Value = -487860.3

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16570.357
$ws.Range("I32").Value = 13434.383
$ws.Range("K32").Value = 13434.383
$ws.Range("M32").Value = -13147.383

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1787.85
$ws.Range("I20").Value = 1816.2727
$ws.Range("K20").Value = 1816.2727
$ws.Range("M20").Value = -1569.2727
$ws.Range("H75").Value = 36471.93
$ws.Range("I75").Value = 7785.3335
$ws.Range("J75").Value = 57986.875
$ws.Range("K75").Value = 7785.3335
$ws.Range("L75").Value = 57986.875
$ws.Range("M75").Value = -6849.3335
$ws.Range("N75").Value = -59858.875
$ws.Range("H78").Value = 36471.93
$ws.Range("I78").Value = 7785.3335
$ws.Range("J78").Value = 57986.875
$ws.Range("K78").Value = 23356.0005
$ws.Range("L78").Value = 173960.625
$ws.Range("M78").Value = -18676.0005
$ws.Range("N78").Value = -183320.625

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 72697
$ws.Range("J88").Value = 72697
$ws.Range("L88").Value = 72697
$ws.Range("N88").Value = -73509
$ws.Range("H91").Value = 72697
$ws.Range("J91").Value = 72697
$ws.Range("L91").Value = 72697
$ws.Range("N91").Value = -75505
$ws.Range("H99").Value = 4222
$ws.Range("I99").Value = 3675.9092
$ws.Range("J99").Value = 5723.75
$ws.Range("K99").Value = 3675.9092
$ws.Range("L99").Value = 5723.75
$ws.Range("M99").Value = -2177.9092
$ws.Range("N99").Value = -8719.75
$ws.Range("H126").Value = 4222
$ws.Range("I126").Value = 3675.9092
$ws.Range("J126").Value = 5723.75
$ws.Range("K126").Value = 11027.7276
$ws.Range("L126").Value = 17171.25
$ws.Range("M126").Value = -8557.7276
$ws.Range("N126").Value = -22111.25
$ws.Range("H141").Value = 242105
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 242105
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 242105
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -252465

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 426.6154
$ws.Range("J2").Value = 453.91666
$ws.Range("L2").Value = 2723.49996
$ws.Range("N2").Value = -2949.49996
$ws.Range("H7").Value = 2056.2
$ws.Range("I7").Value = 1319.75
$ws.Range("K7").Value = 3959.25
$ws.Range("M7").Value = -3847.25
$ws.Range("H34").Value = 391.83334
$ws.Range("I34").Value = 391.83334
$ws.Range("K34").Value = 1175.50002
$ws.Range("M34").Value = -1091.50002
$ws.Range("H37").Value = 47740
$ws.Range("J37").Value = 47740
$ws.Range("L37").Value = 143220
$ws.Range("N37").Value = -143444
$ws.Range("H55").Value = 113447.22
$ws.Range("I55").Value = 506.25
$ws.Range("J55").Value = 203800
$ws.Range("K55").Value = 1518.75
$ws.Range("L55").Value = 611400
$ws.Range("M55").Value = -1341.75
$ws.Range("N55").Value = -611754
$ws.Range("H62").Value = 7000
$ws.Range("I62").Value = 7000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 21000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -20314
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 7000
$ws.Range("I65").Value = 7000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 63000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -59568
$ws.Range("N65").ClearContents()
$ws.Range("H113").Value = 1999.6666
$ws.Range("J113").Value = 1999
$ws.Range("L113").Value = 5997
$ws.Range("N113").Value = -10337

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 20015120
$ws.Range("I70").Value = 25004350
$ws.Range("K70").Value = 25004350
$ws.Range("M70").Value = -25004080
$ws.Range("H73").Value = 20015120
$ws.Range("I73").Value = 25004350
$ws.Range("K73").Value = 25004350
$ws.Range("M73").Value = -25003414
$ws.Range("H80").Value = 2722031.5
$ws.Range("I80").Value = 4237180
$ws.Range("J80").Value = 701833.3
$ws.Range("K80").Value = 4237180
$ws.Range("L80").Value = 701833.3
$ws.Range("M80").Value = -4236182
$ws.Range("N80").Value = -703829.3
$ws.Range("H83").Value = 2722031.5
$ws.Range("I83").Value = 4237180
$ws.Range("J83").Value = 701833.3
$ws.Range("K83").Value = 21185900
$ws.Range("L83").Value = 3509166.5
$ws.Range("M83").Value = -21180908
$ws.Range("N83").Value = -3519150.5
$ws.Range("H104").Value = 39999
$ws.Range("J104").Value = 39999
$ws.Range("L104").Value = 39999
$ws.Range("N104").Value = -46987
$ws.Range("H122").Value = 332530.84
$ws.Range("I122").Value = 373347.25
$ws.Range("K122").Value = 1120041.75
$ws.Range("M122").Value = -1117591.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11433.777
$ws.Range("I7").Value = 10984.333
$ws.Range("J7").Value = 11658.5
$ws.Range("K7").Value = 10984.333
$ws.Range("L7").Value = 11658.5
$ws.Range("M7").Value = -10872.333
$ws.Range("N7").Value = -11882.5
$ws.Range("H40").Value = 8311.458000000001
$ws.Range("I40").Value = 5252.3335
$ws.Range("J40").Value = 11370.583
$ws.Range("K40").Value = 5252.3335
$ws.Range("L40").Value = 11370.583
$ws.Range("M40").Value = -5116.3335
$ws.Range("N40").Value = -11642.583
$ws.Range("H122").Value = 8640.066000000001
$ws.Range("J122").Value = 7937
$ws.Range("L122").Value = 23811
$ws.Range("N122").Value = -28711
$ws.Range("H126").Value = 11433.777
$ws.Range("I126").Value = 10984.333
$ws.Range("J126").Value = 11658.5
$ws.Range("K126").Value = 32952.999
$ws.Range("L126").Value = 34975.5
$ws.Range("M126").Value = -30482.999
$ws.Range("N126").Value = -39915.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10160.814
$ws.Range("I62").Value = 13755.444
$ws.Range("K62").Value = 13755.444
$ws.Range("M62").Value = -13131.444
$ws.Range("H65").Value = 10160.814
$ws.Range("I65").Value = 13755.444
$ws.Range("K65").Value = 68777.22
$ws.Range("M65").Value = -65657.22
$ws.Range("H100").Value = 4400.5293
$ws.Range("I100").Value = 5552.8335
$ws.Range("J100").Value = 1635
$ws.Range("K100").Value = 11105.667
$ws.Range("L100").Value = 3270
$ws.Range("M100").Value = -10564.667
$ws.Range("N100").Value = -4352
$ws.Range("H113").Value = 766.6129
$ws.Range("I113").Value = 261.11765
$ws.Range("J113").Value = 1380.4286
$ws.Range("K113").Value = 783.3529500000001
$ws.Range("L113").Value = 4141.2858
$ws.Range("M113").Value = 1386.64705
$ws.Range("N113").Value = -8481.2858
$ws.Range("H132").Value = 32615078
$ws.Range("I132").Value = 40006100
$ws.Range("J132").Value = 1819157.5
$ws.Range("K132").Value = 120018300
$ws.Range("L132").Value = 5457472.5
$ws.Range("M132").Value = -120015770
$ws.Range("N132").Value = -5462532.5
